# Updated cryptos list refresh (prices + 1h volume deltas), including the
# ShibaInu/Avalanche row-order swap, matching the upstream scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.548.88'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.646.00'
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '596.30'
$ws.Range('E5').Value = '  -0.78%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '156.13'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.629'
$ws.Range('E8').Value = '  +2.07%  '
$ws.Range('E9').Value = '  +2.93%  '
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('E11').Value = '  -1.41%  '
$ws.Range('E12').Value = '  +1.03%  '
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000198'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.63'
$ws.Range('E14').Value = '  -3.23%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.121.12'
$ws.Range('E15').Value = '  -1.06%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.393.47'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.629.29'
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.72'
$ws.Range('E19').Value = '  -2.23%  '
$ws.Range('E20').Value = '  -2.22%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '348.34'
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.80'
$ws.Range('E23').Value = '  -1.63%  '
$ws.Range('E24').Value = '  +1.24%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.59'
$ws.Range('E25').Value = '  -2.58%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.68'
$ws.Range('E26').Value = '  +2.81%  '
$ws.Range('E27').Value = '  -1.72%  '
$ws.Range('E28').Value = '  -2.59%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  -3.52%  '
$ws.Range('E31').Value = '  -1.47%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '526.18'
$ws.Range('E32').Value = '  -3.15%  '
$ws.Range('E33').Value = '  -1.68%  '
$ws.Range('E34').Value = '  -3.29%  '
$ws.Range('E35').Value = '  -1.55%  '
$ws.Range('E36').Value = '  -1.17%  '
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('E39').Value = '  -1.86%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '155.29'
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '160.77'
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('E43').Value = '  -1.04%  '
$ws.Range('E44').Value = '  -2.14%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.26'
$ws.Range('E45').Value = '  -1.88%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '22.47'
$ws.Range('E46').Value = '  -3.49%  '
$ws.Range('E47').Value = '  -1.89%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0253'
$ws.Range('E48').Value = '  -2.66%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0994'
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0₆0250'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '19.67'
$ws.Range('E51').Value = '  -2.37%  '
